# Insert two new rows at 702-703 (shifts old rows 702.. down to 704..)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("702:703").Insert()

# Fill new row 702 - Coliflor, Primera, market date 2022-06-02 (serial 44714)
$ws.Range("A702").Value2 = 6
$ws.Range("B702").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C702").Value2 = "Metropolitana"
$ws.Range("D702").Value2 = 44714
$ws.Range("E702").Value2 = 13
$ws.Range("F702").Value2 = 100112008
$ws.Range("G702").Value2 = "Coliflor"
$ws.Range("H702").Value2 = "Sin especificar"
$ws.Range("I702").Value2 = "Primera"
$ws.Range("J702").Value2 = 11400
$ws.Range("K702").Value2 = 700
$ws.Range("L702").Value2 = 900
$ws.Range("M702").Value2 = 801
$ws.Range("N702").Value2 = "$/unidad"
$ws.Range("O702").Value2 = "Región Metropolitana"
$ws.Range("P702").Value2 = 801
$ws.Range("Q702").Value2 = 1
$ws.Range("R702").Value2 = "Hortaliza"

# Fill new row 703 - Coliflor, Segunda, market date 2022-06-02 (serial 44714)
$ws.Range("A703").Value2 = 6
$ws.Range("B703").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C703").Value2 = "Metropolitana"
$ws.Range("D703").Value2 = 44714
$ws.Range("E703").Value2 = 13
$ws.Range("F703").Value2 = 100112008
$ws.Range("G703").Value2 = "Coliflor"
$ws.Range("H703").Value2 = "Sin especificar"
$ws.Range("I703").Value2 = "Segunda"
$ws.Range("J703").Value2 = 5500
$ws.Range("K703").Value2 = 500
$ws.Range("L703").Value2 = 600
$ws.Range("M703").Value2 = 558
$ws.Range("N703").Value2 = "$/unidad"
$ws.Range("O703").Value2 = "Región Metropolitana"
$ws.Range("P703").Value2 = 558
$ws.Range("Q703").Value2 = 1
$ws.Range("R703").Value2 = "Hortaliza"
